$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New ConditionType (C) and ITI (D) values for trials 1-16 (rows 2-17)
$data = @(
    @{ C = 3; D = 7 },
    @{ C = 2; D = 6 },
    @{ C = 4; D = 6 },
    @{ C = 1; D = 6 },
    @{ C = 2; D = 8 },
    @{ C = 2; D = 8 },
    @{ C = 3; D = 8 },
    @{ C = 3; D = 9 },
    @{ C = 4; D = 6 },
    @{ C = 4; D = 6 },
    @{ C = 2; D = 7 },
    @{ C = 1; D = 6 },
    @{ C = 4; D = 9 },
    @{ C = 1; D = 7 },
    @{ C = 1; D = 7 },
    @{ C = 3; D = 6 }
)

# Add ITI header
$ws.Cells.Item(1, 4).Value = "ITI"

# Update ConditionType column and fill in ITI column for rows 2-17
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i].C
    $ws.Cells.Item($row, 4).Value = $data[$i].D
}

# Remove the old trailing rows (Trial 17, 18, 19), previously rows 18-20
$ws.Range("A18:A20").EntireRow.Delete()

$ws.Range("C18").Select()
